$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Unmerge the two "employee" columns that used to span 3 report rows
#    (B9:B11 / C9:C11) so each becomes a single, independent cell.
# ---------------------------------------------------------------------------
$ws.Range("B9:B11").UnMerge() | Out-Null
$ws.Range("C9:C11").UnMerge() | Out-Null

# ---------------------------------------------------------------------------
# 2. Drop the now-obsolete sample rows 10-15 (the template used to ship with
#    3 hard-coded example employees/stores plus a few blank rows below).
#    Only row 9 remains as the single templated data row.
# ---------------------------------------------------------------------------
$ws.Rows("10:15").Delete() | Out-Null

# ---------------------------------------------------------------------------
# 3. Re-style B9/C9 like the rest of the data row: plain thin box border
#    (same as A9), vertically centered, general horizontal alignment, no
#    wrap - instead of the old "merged cell" box border + centered + wrap
#    formatting. Copy the box-border/font from A9 first (keeps the existing
#    border definition instead of synthesizing a new one), then fix up the
#    alignment to vertical-center only.
# ---------------------------------------------------------------------------
$ws.Range("A9").Copy() | Out-Null
$ws.Range("B9:C9").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("B9:C9").HorizontalAlignment = 1          # xlHAlignGeneral
$ws.Range("B9:C9").VerticalAlignment = -4108        # xlVAlignCenter
$ws.Range("B9:C9").WrapText = $false

# ---------------------------------------------------------------------------
# 4. Re-point every label / placeholder to the new template vocabulary.
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "{{Start}}"
$ws.Range("E5").Value = "{{End}}"

$ws.Range("E7").Value = "Đại lý"

$ws.Range("A8").Value = "{{MonitorStoreImages.OrganizationName}}"

$ws.Range("A9").Value = "{{MonitorStoreImages.SaleEmployees.StoreCheckings.STT}}"
$ws.Range("B9").Value = "{{MonitorStoreImages.SaleEmployees.Username}}"
$ws.Range("C9").Value = "{{MonitorStoreImages.SaleEmployees.DisplayName}}"
$ws.Range("D9").Value = "{{MonitorStoreImages.SaleEmployees.StoreCheckings.DateDisplay}}"
$ws.Range("E9").Value = "{{MonitorStoreImages.SaleEmployees.StoreCheckings.StoreName}}"
$ws.Range("F9").Value = "{{MonitorStoreImages.SaleEmployees.StoreCheckings.ImageCounter}}"

# ---------------------------------------------------------------------------
# 5. Selection cosmetics to match the saved file (active cell moved up once
#    the extra rows disappeared).
# ---------------------------------------------------------------------------
$ws.Range("E5").Select() | Out-Null
